$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.085.66'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '1.564.14'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  +0.56%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.490'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  +0.49%  '
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0861'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').Value = '1.788.28'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '1.560.56'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('D16').Value = '27.082.66'
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.99'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('D18').Value = '0.0₃0701'
$ws.Range('E18').Value = '  -0.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '214.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.43%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.82'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.01'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('E30').Value = '  +4.25%  '
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('E33').Value = '  +1.86%  '
$ws.Range('D34').Value = '1.440.14'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.80'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.805'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.34%  '
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('D47').Value = '1.702.06'
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.96'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('E49').Value = '  +3.98%  '
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0957'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.33%  '
